# Insert a new weekly price record at row 392 ("Hortaliza, Macroferia
# Regional de Talca - Acelga"). All existing rows from 392 downward shift
# down by one (handled automatically by the row insert), and the new row
# is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 392 (and everything below it) down by one row.
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row 392 with the new record.
$ws.Range("A392").Value = 5
$ws.Range("B392").Value = "Macroferia Regional de Talca"
$ws.Range("C392").Value = "Maule"
$ws.Range("D392").Value = 45215
$ws.Range("E392").Value = 7
$ws.Range("F392").Value = 100112009
$ws.Range("G392").Value = "Acelga"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 500
$ws.Range("K392").Value = 1800
$ws.Range("L392").Value = 1800
$ws.Range("M392").Value = 1800
$ws.Range("N392").Value = '$/docena de atados (4 kilos)'
$ws.Range("O392").Value = "Región del Maule"
$ws.Range("P392").Value = 450
$ws.Range("Q392").Value = 4
$ws.Range("R392").Value = "Hortaliza"
